$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.873.77'
$ws.Range("E2").Value = '  -2.35%  '

$ws.Range("D3").Value = '3.945.26'
$ws.Range("E3").Value = '  -3.20%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '539.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.67%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.24'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").Value = '3.940.24'
$ws.Range("E7").Value = '  -3.13%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.683'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -7.83%  '

$ws.Range("E9").Value = '  +0.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.735'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.31%  '

$ws.Range("E11").Value = '  -6.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.76'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +15.63%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000316'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.53'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.77%  '

$ws.Range("D15").Value = '4.575.81'
$ws.Range("E15").Value = '  -2.71%  '

$ws.Range("D16").Value = '3.938.73'
$ws.Range("E16").Value = '  -3.11%  '

$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.33'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.53%  '

$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.87%  '

$ws.Range("E19").Value = '  -1.28%  '

$ws.Range("E20").Value = '  -4.84%  '

$ws.Range("D21").Value = '70.799.14'
$ws.Range("E21").Value = '  -2.26%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '422.52'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.65%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.59'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.95%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '96.67'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -7.95%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.20'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.61%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.17'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.41%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.49'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.66'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.20%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.69'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +16.03%  '

$ws.Range("B30").Value = 'LEO'
$ws.Range("C30").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.84'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.64%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '36.33'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.77'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +13.17%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '50.60'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +16.25%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '694.09'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.73%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.32'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.22%  '

$ws.Range("E36").Value = '  -1.38%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '64.47'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.38%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.434'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.52%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.46'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.32%  '

$ws.Range("D40").Value = '0.0₃0816'
$ws.Range("E40").Value = '  -5.60%  '

$ws.Range("E41").Value = '  -1.98%  '

$ws.Range("E42").Value = '  +0.11%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.14%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.20'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.47%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0479'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.98%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.72'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.29%  '

$ws.Range("E47").Value = '  -8.33%  '

$ws.Range("E48").Value = '  +6.12%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.39'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.96%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.97'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.63%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000273'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.98%  '
